$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed Coin/Link/Price/Volume cells to reflect the refreshed
# cryptocurrency data snapshot (GitHub Actions scheduled update).

$ws.Range("D2").Value = "'30.409.45"
$ws.Range("E2").Value = "  -0.23%  "
$ws.Range("D3").Value = "'2.072.19"
$ws.Range("E3").Value = "  +3.71%  "
$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("D5").Value = "'328.80"
$ws.Range("E5").Value = "  +1.32%  "
$ws.Range("D6").Value = "'1.005"
$ws.Range("E6").Value = "  +0.39%  "
$ws.Range("D7").Value = "'0.5198"
$ws.Range("E7").Value = "  +1.79%  "
$ws.Range("D8").Value = "'0.4348"
$ws.Range("E8").Value = "  +4.98%  "
$ws.Range("D9").Value = "'0.08649"
$ws.Range("E9").Value = "  -0.72%  "
$ws.Range("D10").Value = "'45.88"
$ws.Range("E10").Value = "  +6.70%  "
$ws.Range("D11").Value = "'1.151"
$ws.Range("E11").Value = "  +1.67%  "
$ws.Range("D12").Value = "'24.16"
$ws.Range("E12").Value = "  -2.16%  "
$ws.Range("D13").Value = "'2.080.59"
$ws.Range("E13").Value = "  +4.27%  "
$ws.Range("D14").Value = "'6.615"
$ws.Range("E14").Value = "  +0.60%  "
$ws.Range("D15").Value = "'7.684"
$ws.Range("E15").Value = "  +3.39%  "
$ws.Range("D16").Value = "'1.006"
$ws.Range("E16").Value = "  +0.40%  "
$ws.Range("D17").Value = "'95.13"
$ws.Range("E17").Value = "  +0.90%  "
$ws.Range("E18").Value = "  -0.36%  "
$ws.Range("E19").Value = "  +1.61%  "
$ws.Range("D20").Value = "'18.70"
$ws.Range("E20").Value = "  -1.34%  "
$ws.Range("D21").Value = "'1.005"
$ws.Range("E21").Value = "  +0.44%  "
$ws.Range("D22").Value = "'6.229"
$ws.Range("E22").Value = "  +1.01%  "
$ws.Range("D23").Value = "'30.445.63"
$ws.Range("E23").Value = "  -0.25%  "
$ws.Range("E24").Value = "  +3.46%  "
$ws.Range("D25").Value = "'2.311"
$ws.Range("E25").Value = "  +3.82%  "
$ws.Range("D26").Value = "'2.325.91"
$ws.Range("E26").Value = "  +4.43%  "
$ws.Range("D27").Value = "'22.11"
$ws.Range("E27").Value = "  -1.17%  "
$ws.Range("D28").Value = "'161.55"
$ws.Range("E28").Value = "  -1.14%  "
$ws.Range("D29").Value = "'2.511"
$ws.Range("E29").Value = "  +4.61%  "
$ws.Range("D30").Value = "'130.46"
$ws.Range("E30").Value = "  -0.76%  "
$ws.Range("D31").Value = "'1.173"
$ws.Range("E31").Value = "  +3.23%  "
$ws.Range("D32").Value = "'0.1068"
$ws.Range("E32").Value = "  +1.60%  "
$ws.Range("D33").Value = "'6.029"
$ws.Range("E33").Value = "  -0.77%  "
$ws.Range("D34").Value = "'3.818"
$ws.Range("E34").Value = "  -0.93%  "
$ws.Range("D35").Value = "'1.494"
$ws.Range("E35").Value = "  +12.08%  "
$ws.Range("D36").Value = "'0.02556"
$ws.Range("E36").Value = "  +1.60%  "
$ws.Range("D37").Value = "'9.597"
$ws.Range("D38").Value = "'5.435"
$ws.Range("E38").Value = "  +0.10%  "
$ws.Range("D39").Value = "'0.06589"
$ws.Range("E39").Value = "  -0.17%  "
$ws.Range("D40").Value = "'12.41"
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("E41").Value = "  +1.35%  "
$ws.Range("D42").Value = "'0.6698"
$ws.Range("E42").Value = "  +1.21%  "
$ws.Range("D43").Value = "'1.233"
$ws.Range("E43").Value = "  -0.33%  "
$ws.Range("E44").Value = "  +0.38%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "'0.6293"
$ws.Range("E45").Value = "  +2.11%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'13.80"
$ws.Range("E46").Value = "  +1.66%  "
$ws.Range("D47").Value = "'2.188"
$ws.Range("E47").Value = "  -0.47%  "
$ws.Range("E48").Value = "  -1.17%  "
$ws.Range("D49").Value = "'1.228"
$ws.Range("E49").Value = "  -3.12%  "
$ws.Range("D50").Value = "'81.46"
$ws.Range("E50").Value = "  +1.73%  "
$ws.Range("D51").Value = "'1.178"
$ws.Range("E51").Value = "  +6.37%  "
